$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder category rows: Transport, Food, Bills (row2..row4)
$ws.Range("A2").Value = "Transport"
$ws.Range("B2").Value = 156

$ws.Range("A3").Value = "Food"
$ws.Range("B3").Value = 143

$ws.Range("A4").Value = "Bills"
$ws.Range("B4").Value = 900

# Update Total Spent and Remaining
$ws.Range("B7").Value = 1199
$ws.Range("B8").Value = 548801
